$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "[1, 0, 0, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'RegulationViolation']"

$ws.Range("D24").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E24").Value = "['HardwareFault']"

$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E25").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

$ws.Range("D28").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E28").Value = "['SoftwareFault']"

$ws.Range("D29").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E29").Value = "['SoftwareFault']"

$ws.Range("D39").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D53").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal']"

$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 1, 0]"
$ws.Range("E54").Value = "['CommunicationIssue']"

$ws.Range("D56").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "[]"

$ws.Range("D69").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal']"

$ws.Range("D81").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E81").Value = "['Normal', 'HardwareFault']"

$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 1]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment', 'SoftwareFault']"

$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault']"

$ws.Range("D118").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E118").Value = "['Normal']"
